$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (A: numeric id, B: upn string)
$data = @(
    @(0, "M931325212046"),
    @(1, "X931235209022"),
    @(2, "C931100609010"),
    @(4, "M931321110016")
)

$startRow = 3
$endRow = $startRow + $data.Count - 1

# Clone the formatting used by the existing A2 cell (bold/border/centered)
# onto the new A-column cells before writing their values.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A$startRow`:A$endRow").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
